# Add a new sheet ("Sheet1") after the existing "Sayfa1" sheet and populate it
# with a header row (copied from Sayfa1) plus one new component row, then make
# it the active sheet - matching the source-controlled diff.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sayfa1")

# Select the header row on the original sheet (mirrors the recorded selection
# state left behind on Sayfa1 after copying the header row: A1:F1).
$ws1.Range("A1:F1").Select()

# Add the new worksheet right after Sayfa1, becoming tab 2.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet1"

# Header row, identical to Sayfa1's header row.
$ws2.Range("A1").Value = "Digikey Part No"
$ws2.Range("B1").Value = "Description"
$ws2.Range("C1").Value = "Digikey Link"
$ws2.Range("D1").Value = "Adet"
$ws2.Range("E1").Value = "Price"
$ws2.Range("F1").Value = "Total Price"
$ws2.Range("A1:F1").Font.Bold = $true

# New component row.
$ws2.Range("A2").Value = "P.024AUCT-ND"
$ws2.Range("B2").Value = "RES 0.024 OHM 1% 1W 1206"
$ws2.Range("C2").Value = "https://www.digikey.com/product-detail/en/panasonic-electronic-components/ERJ-8BWFR024V/P.024AUCT-ND/1711691"
$ws2.Range("D2").Value = 1
$ws2.Range("E2").NumberFormat = "@"
$ws2.Range("E2").Value = "0.74"

$ws2.Columns.Item(1).AutoFit() | Out-Null
$ws2.Columns.Item(2).AutoFit() | Out-Null
$ws2.Columns.Item(3).AutoFit() | Out-Null
$ws2.Columns.Item(4).AutoFit() | Out-Null
$ws2.Columns.Item(5).AutoFit() | Out-Null
$ws2.Columns.Item(6).AutoFit() | Out-Null

$ws2.Range("F2").Select()
$ws2.Activate()
